$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value()
$text = $text.Replace("✅ 1000 Bs = 1.89 = 6814.93 pesos", "✅ 1000 Bs = 1.87 = 6735.96 pesos")
$text = $text.Replace("✅ 6814.93 pesos = 1.88 = 955.33 Bs", "✅ 6735.96 pesos = 1.86 = 955.28 Bs")
$cellA1.Value = $text

# --- Update the rate figures on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 534
$wsTasas.Range("O10").Value = 3597
$wsTasas.Range("N12").Value = 3618
$wsTasas.Range("O12").Value = 513.1
